$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "Tìm hiểu đề tài, công nghệ sử dụng, vẽ sơ đồ uscase tổng quát"
$ws.Range("D3").Value = "mô tả usecase, làm sơ đồ erd"
$ws.Range("D2").WrapText = $true
$ws.Range("D2").Select()
